$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "Lista de asistencia hola"

# Update row 2 values (attendance record replaced by a new one)
$ws.Range("A2").Value = "Samuel Hiram Castro Martinez"

# Keep "No. de control" as text (it was text before too)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20212390"

$ws.Range("C2").Value = "ING SISTEMAS"

# Date string stored as plain text (matches source format)
$ws.Range("D2").Value = "04/07/2024 14:42:46"

# Numeric count of attendances
$ws.Range("E2").Value = 1

$ws.Range("F2").Value = "Samuel Castro Martinez"
